# Clean System Level v0.1
# Applies the data + selection changes captured in the commit diff:
#  - technology_fuel_pairs: BF-BOF fuel shares for Coal/COG updated, sheet
#    becomes the active tab with a new selection.
#  - technology_material_pairs: no data change; it simply stops being the
#    active tab (handled implicitly once another sheet is activated).
#  - fuel_cost: Coal (row 2) ramped up to 5 across all years, Bio (row 6)
#    zeroed out across all years.
#  - emission: the "global" emission cap row gets a step-down formula chain
#    (halving every 5 years) through 2044, then drops to 0 from 2045 on.

$wb = $excel.ActiveWorkbook

# --- technology_fuel_pairs ---
$wsTfp = $wb.Worksheets.Item("technology_fuel_pairs")
$wsTfp.Range("C2").Value = 0.5
$wsTfp.Range("C3").Value = 0.7

# --- fuel_cost ---
$wsFuelCost = $wb.Worksheets.Item("fuel_cost")
$wsFuelCost.Activate() | Out-Null
$wsFuelCost.Range("B2:AA2").Value = 5
$wsFuelCost.Range("B6:AA6").Value = 0
$wsFuelCost.Range("B6:AA6").Select() | Out-Null

# --- emission ---
$wsEmission = $wb.Worksheets.Item("emission")
$wsEmission.Activate() | Out-Null
$wsEmission.Range("C2").Formula = "=B2*0.5"
$wsEmission.Range("D2").Formula = "=C2"
$wsEmission.Range("E2:F2").Formula = "=D2"
$wsEmission.Range("G2").Formula = "=F2*0.5"
$wsEmission.Range("H2").Formula = "=G2"
$wsEmission.Range("I2:K2").Formula = "=H2"
$wsEmission.Range("L2").Formula = "=K2*0.5"
$wsEmission.Range("M2").Formula = "=L2"
$wsEmission.Range("N2:U2").Formula = "=M2"
$wsEmission.Range("V2:AA2").Value = 0
$wsEmission.Range("V2:AA2").Select() | Out-Null

# --- technology_fuel_pairs becomes (and stays) the active tab/selection ---
$wsTfp.Activate() | Out-Null
$wsTfp.Range("S10").Select() | Out-Null
